# The deck's single Design ("Integral") is being switched to the stock
# "Office Theme" colour palette. All structural parts of the theme
# (font scheme / format scheme) are already identical between the two
# themes in this file, so only the 12 theme colours need to change.
#
# Office Theme palette, in PowerPoint's dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink order, expressed as COM RGB() long values (R + G*256 + B*65536):
#   dk1      000000 -> 0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A -> 6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 -> 3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 -> 49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 -> 4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 -> 7491477

$p = $ppt.ActivePresentation

$officeThemeRGB = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$s = $p.Slides.Item(1)
$themeColors = $s.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
